$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''28.238.14'
$ws.Range("E2").Value = '  +3.26%  '
$ws.Range("D3").Value = '''1.923.52'
$ws.Range("E3").Value = '  +3.12%  '
$ws.Range("E4").Value = '  -1.52%  '
$ws.Range("D5").Value = '''316.43'
$ws.Range("E5").Value = '  +0.71%  '
$ws.Range("E6").Value = '  -1.40%  '
$ws.Range("D7").Value = '''0.4852'
$ws.Range("E7").Value = '  +0.88%  '
$ws.Range("D8").Value = '''0.3851'
$ws.Range("E8").Value = '  +3.12%  '
$ws.Range("D9").Value = '''0.07405'
$ws.Range("E9").Value = '  -0.05%  '
$ws.Range("D10").Value = '''0.9463'
$ws.Range("E10").Value = '  +0.84%  '
$ws.Range("D11").Value = '''20.92'
$ws.Range("E11").Value = '  +0.63%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '''1.995.13'
$ws.Range("E12").Value = '  +6.57%  '
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").Value = '''0.07797'
$ws.Range("E13").Value = '  -1.14%  '
$ws.Range("D14").Value = '''5.538'
$ws.Range("E14").Value = '  +1.77%  '
$ws.Range("D15").Value = '''6.679'
$ws.Range("E15").Value = '  +1.83%  '
$ws.Range("D16").Value = '''92.06'
$ws.Range("E16").Value = '  +1.92%  '
$ws.Range("D17").Value = '''1.008'
$ws.Range("E17").Value = '  -1.43%  '
$ws.Range("D18").Value = '''0.000008876'
$ws.Range("E18").Value = '  +0.77%  '
$ws.Range("D19").Value = '''1.007'
$ws.Range("E19").Value = '  -1.24%  '
$ws.Range("D20").Value = '''28.258.96'
$ws.Range("E20").Value = '  +3.15%  '
$ws.Range("E21").Value = '  +1.46%  '
$ws.Range("D22").Value = '''5.173'
$ws.Range("E22").Value = '  +0.73%  '
$ws.Range("D23").Value = '''2.184.94'
$ws.Range("E23").Value = '  +3.41%  '
$ws.Range("E24").Value = '  +2.49%  '
$ws.Range("E25").Value = '  -1.82%  '
$ws.Range("D26").Value = '''156.34'
$ws.Range("E26").Value = '  +1.36%  '
$ws.Range("D27").Value = '''18.66'
$ws.Range("E27").Value = '  +0.64%  '
$ws.Range("D28").Value = '''2.114'
$ws.Range("E28").Value = '  +4.95%  '
$ws.Range("D29").Value = '''117.09'
$ws.Range("E29").Value = '  +0.92%  '
$ws.Range("D30").Value = '''5.035'
$ws.Range("E30").Value = '  +0.76%  '
$ws.Range("E31").Value = '  -0.28%  '
$ws.Range("E32").Value = '  +0.31%  '
$ws.Range("D33").Value = '''1.256'
$ws.Range("E33").Value = '  +4.91%  '
$ws.Range("D34").Value = '''0.7762'
$ws.Range("E34").Value = '  +4.24%  '
$ws.Range("D35").Value = '''4.690'
$ws.Range("E35").Value = '  +2.76%  '
$ws.Range("D36").Value = '''2.771'
$ws.Range("E36").Value = '  +2.87%  '
$ws.Range("D37").Value = '''0.02057'
$ws.Range("E37").Value = '  +0.07%  '
$ws.Range("D38").Value = '''1.130'
$ws.Range("E38").Value = '  +0.46%  '
$ws.Range("B39").Value = 'TheSandbox'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D39").Value = '''0.5604'
$ws.Range("E39").Value = '  +4.40%  '
$ws.Range("B40").Value = 'Hedera'
$ws.Range("C40").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D40").Value = '''0.05372'
$ws.Range("E40").Value = '  +1.32%  '
$ws.Range("D41").Value = '''3.043'
$ws.Range("E41").Value = '  +1.39%  '
$ws.Range("D42").Value = '''7.109'
$ws.Range("E42").Value = '  -0.20%  '
$ws.Range("D43").Value = '''8.584'
$ws.Range("E43").Value = '  +2.14%  '
$ws.Range("D44").Value = '''0.1540'
$ws.Range("D45").Value = '''0.4923'
$ws.Range("E45").Value = '  +1.76%  '
$ws.Range("D46").Value = '''10.77'
$ws.Range("E46").Value = '  +0.79%  '
$ws.Range("D47").Value = '''107.26'
$ws.Range("E47").Value = '  +3.77%  '
$ws.Range("E48").Value = '  -1.50%  '
$ws.Range("D49").Value = '''1.679'
$ws.Range("E49").Value = '  +0.52%  '
$ws.Range("D50").Value = '''69.68'
$ws.Range("E50").Value = '  +4.11%  '
$ws.Range("D51").Value = '''0.06156'
$ws.Range("E51").Value = '  +0.92%  '
